$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "NA" placeholders in September/October/November/December 2020 (J6:M6)
# are removed now that the real sales figures/charts are in place.
$ws.Range("J6:M6").ClearContents()

# Leave the selection where the author left it after the edit.
$ws.Range("M6").Select()
